$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44498
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 1600
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = 1600
$ws.Range("P2").Value = 1600

$ws.Range("D3").Value = 44475
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 2000
$ws.Range("P3").Value = 2000

$ws.Range("D4").Value = 44161
$ws.Range("H4").Value = "Verde"
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 1700
$ws.Range("L4").Value = 1700
$ws.Range("M4").Value = 1700
$ws.Range("P4").Value = 1700

$ws.Range("D5").Value = 44162
$ws.Range("H5").Value = "Verde"
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 1600
$ws.Range("L5").Value = 1600
$ws.Range("M5").Value = 1600
$ws.Range("P5").Value = 1600

$ws.Range("D6").Value = 44504
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("J6").Value = 180
$ws.Range("K6").Value = 1600
$ws.Range("L6").Value = 1600
$ws.Range("M6").Value = 1600
$ws.Range("P6").Value = 1600

$ws.Range("D7").Value = 44516
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("J7").Value = 360
$ws.Range("K7").Value = 1600
$ws.Range("L7").Value = 1600
$ws.Range("M7").Value = 1600
$ws.Range("P7").Value = 1600

$ws.Range("D8").Value = 44488
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 1700
$ws.Range("L8").Value = 1800
$ws.Range("M8").Value = 1750
$ws.Range("P8").Value = 1750

$ws.Range("D9").Value = 44526
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("J9").Value = 600
$ws.Range("K9").Value = 1700
$ws.Range("L9").Value = 1700
$ws.Range("M9").Value = 1700
$ws.Range("P9").Value = 1700

$ws.Range("D10").Value = 44523
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("J10").Value = 520
$ws.Range("K10").Value = 1800
$ws.Range("L10").Value = 1800
$ws.Range("M10").Value = 1800
$ws.Range("P10").Value = 1800

$ws.Range("D11").Value = 44481
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 1700
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1850
$ws.Range("P11").Value = 1850

$ws.Range("D12").Value = 44482
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("J12").Value = 72
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = 2000
$ws.Range("P12").Value = 2000

$ws.Range("D13").Value = 44169
$ws.Range("H13").Value = "Verde"
$ws.Range("J13").Value = 600
$ws.Range("K13").Value = 1600
$ws.Range("L13").Value = 1600
$ws.Range("M13").Value = 1600
$ws.Range("P13").Value = 1600

$ws.Range("D14").Value = 44495
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("J14").Value = 520
$ws.Range("K14").Value = 1800
$ws.Range("L14").Value = 1800
$ws.Range("M14").Value = 1800
$ws.Range("P14").Value = 1800

$ws.Range("D15").Value = 44511
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("J15").Value = 144
$ws.Range("K15").Value = 1700
$ws.Range("L15").Value = 1700
$ws.Range("M15").Value = 1700
$ws.Range("P15").Value = 1700

$ws.Range("D16").Value = 44517
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 1600
$ws.Range("L16").Value = 1600
$ws.Range("M16").Value = 1600
$ws.Range("P16").Value = 1600

$ws.Range("D17").Value = 44179
$ws.Range("H17").Value = "Verde"
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 1600
$ws.Range("L17").Value = 1600
$ws.Range("M17").Value = 1600
$ws.Range("P17").Value = 1600

$ws.Range("D18").Value = 44487
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 1800
$ws.Range("L18").Value = 1800
$ws.Range("M18").Value = 1800
$ws.Range("P18").Value = 1800

$ws.Range("D19").Value = 44176
$ws.Range("H19").Value = "Verde"
$ws.Range("J19").Value = 700
$ws.Range("K19").Value = 1600
$ws.Range("L19").Value = 1600
$ws.Range("M19").Value = 1600
$ws.Range("P19").Value = 1600

$ws.Range("D20").Value = 44165
$ws.Range("H20").Value = "Verde"
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 1600
$ws.Range("L20").Value = 1600
$ws.Range("M20").Value = 1600
$ws.Range("P20").Value = 1600

$ws.Range("D21").Value = 44530
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("J21").Value = 600
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = 1500
$ws.Range("P21").Value = 1500

$ws.Range("D22").Value = 44484
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 1700
$ws.Range("M22").Value = 1700
$ws.Range("P22").Value = 1700

$ws.Range("D23").Value = 44476
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("J23").Value = 60
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("M23").Value = 2000
$ws.Range("P23").Value = 2000

$ws.Range("D24").Value = 44503
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("J24").Value = 72
$ws.Range("K24").Value = 1600
$ws.Range("L24").Value = 1600
$ws.Range("M24").Value = 1600
$ws.Range("P24").Value = 1600

$ws.Range("D25").Value = 44159
$ws.Range("H25").Value = "Verde"
$ws.Range("J25").Value = 600
$ws.Range("K25").Value = 1600
$ws.Range("L25").Value = 1700
$ws.Range("M25").Value = 1650
$ws.Range("P25").Value = 1650

$ws.Range("D26").Value = 44529
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("J26").Value = 120
$ws.Range("K26").Value = 1700
$ws.Range("L26").Value = 1700
$ws.Range("M26").Value = 1700
$ws.Range("P26").Value = 1700

$ws.Range("D27").Value = 44496
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("J27").Value = 84
$ws.Range("K27").Value = 1800
$ws.Range("L27").Value = 1800
$ws.Range("M27").Value = 1800
$ws.Range("P27").Value = 1800

$ws.Range("D28").Value = 44518
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("J28").Value = 180
$ws.Range("K28").Value = 1600
$ws.Range("L28").Value = 1600
$ws.Range("M28").Value = 1600
$ws.Range("P28").Value = 1600

$ws.Range("D29").Value = 44519
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("J29").Value = 600
$ws.Range("K29").Value = 1600
$ws.Range("L29").Value = 1800
$ws.Range("M29").Value = 1700
$ws.Range("P29").Value = 1700

$ws.Range("D30").Value = 44168
$ws.Range("H30").Value = "Verde"
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 1600
$ws.Range("L30").Value = 1600
$ws.Range("M30").Value = 1600
$ws.Range("P30").Value = 1600

$ws.Range("D31").Value = 44473
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 2000
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = 2000
$ws.Range("P31").Value = 2000

$ws.Range("D32").Value = 44166
$ws.Range("H32").Value = "Verde"
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 1600
$ws.Range("L32").Value = 1600
$ws.Range("M32").Value = 1600
$ws.Range("P32").Value = 1600

$ws.Range("D33").Value = 44522
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("J33").Value = 180
$ws.Range("K33").Value = 1800
$ws.Range("L33").Value = 1800
$ws.Range("M33").Value = 1800
$ws.Range("P33").Value = 1800

$ws.Range("D34").Value = 44497
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("J34").Value = 120
$ws.Range("K34").Value = 1800
$ws.Range("L34").Value = 1800
$ws.Range("M34").Value = 1800
$ws.Range("P34").Value = 1800

$ws.Range("D35").Value = 44525
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("J35").Value = 180
$ws.Range("K35").Value = 1600
$ws.Range("L35").Value = 1600
$ws.Range("M35").Value = 1600
$ws.Range("P35").Value = 1600

$ws.Range("D36").Value = 44494
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("J36").Value = 120
$ws.Range("K36").Value = 1700
$ws.Range("L36").Value = 1700
$ws.Range("M36").Value = 1700
$ws.Range("P36").Value = 1700

$ws.Range("D37").Value = 44509
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("J37").Value = 550
$ws.Range("K37").Value = 1700
$ws.Range("L37").Value = 1700
$ws.Range("M37").Value = 1700
$ws.Range("P37").Value = 1700

$ws.Range("D38").Value = 44491
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 1700
$ws.Range("L38").Value = 1700
$ws.Range("M38").Value = 1700
$ws.Range("P38").Value = 1700

$ws.Range("D39").Value = 44474
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("J39").Value = 250
$ws.Range("K39").Value = 2000
$ws.Range("L39").Value = 2000
$ws.Range("M39").Value = 2000
$ws.Range("P39").Value = 2000

$ws.Range("D40").Value = 44490
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("J40").Value = 72
$ws.Range("K40").Value = 1700
$ws.Range("L40").Value = 1700
$ws.Range("M40").Value = 1700
$ws.Range("P40").Value = 1700
